# "Generate Report for Handback"
#
# The localization-status report is regenerated after a handback occurred:
#   - the "Ready for handoff" status text is replaced everywhere by
#     "Handed back: in sync with en-US"
#   - the per-language sheets (zh-cn / de-de) gain a populated
#     "Latest Target File" (col I, now a hyperlink to the source .md doc)
#     and "Latest Handback File" (col J, the handed-back .xlf) with a
#     real "Latest Handback DateTime" (col K) instead of the 0001-01-01
#     placeholder.
#   - a couple of columns get widened so the new long file names are
#     readable.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$mdUrlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ee668b4c8e2d5be722197661375b377b2e169455/e2e/5daf2757-4342-4441-8a1b-7c039d6e4905.md"
$mdUrlB = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ee668b4c8e2d5be722197661375b377b2e169455/e2e/ffffbe4f7217-6524-4177-bc81-4f7a0057af32.md"
$mdNameA = "5daf2757-4342-4441-8a1b-7c039d6e4905.md"
$mdNameB = "ffffbe4f7217-6524-4177-bc81-4f7a0057af32.md"

# Column widths: ColumnWidth is quantized internally to 1/6ths of a
# character by this engine, so the closest attainable width to the
# original author's 29.9777047293527 is 30 (ColumnWidth 29.1666...);
# 40 is exactly attainable (ColumnWidth 39.1666...).
$wideTarget = 29.166666666666668
$fortyTarget = 39.166666666666664

# ---------------------------------------------------------------------
# Overview sheet: just the status text + column widths change.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Columns.Item(5).ColumnWidth = $wideTarget
$wsOverview.Columns.Item(6).ColumnWidth = $wideTarget

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

# Latest Handback File (col J) + Latest Handback DateTime (col K)
$wsZh.Range("J2").Value = "5daf2757-4342-4441-8a1b-7c039d6e4905.d0ea2da892cb80347cf45bf2890c2371202062f7.zh-cn.xlf"
$wsZh.Range("J3").Value = "5daf2757-4342-4441-8a1b-7c039d6e4905.d0ea2da892cb80347cf45bf2890c2371202062f7.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-25 03:02:46"
$wsZh.Range("K3").Value = "2016-08-25 03:02:46"

# Rebuild the hyperlinks collection so the row-interleaved A2,I2,A3,I3
# order (and therefore the relationship id sequence) matches a fresh
# report generation; this also populates the new "Latest Target File"
# (col I) values/styles.
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdUrlA, $null, $null, $mdNameA)
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrlA, $null, $null, $mdNameA)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $mdUrlB, $null, $null, $mdNameB)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdUrlA, $null, $null, $mdNameA)

$wsZh.Columns.Item(3).ColumnWidth = $wideTarget
$wsZh.Columns.Item(9).ColumnWidth = $fortyTarget
$wsZh.Columns.Item(10).ColumnWidth = $fortyTarget

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# Latest Handback File (col J) + Latest Handback DateTime (col K)
$wsDe.Range("J2").Value = "5daf2757-4342-4441-8a1b-7c039d6e4905.d0ea2da892cb80347cf45bf2890c2371202062f7.de-de.xlf"
$wsDe.Range("J3").Value = "5daf2757-4342-4441-8a1b-7c039d6e4905.d0ea2da892cb80347cf45bf2890c2371202062f7.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-25 03:02:53"
$wsDe.Range("K3").Value = "2016-08-25 03:02:53"

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdUrlA, $null, $null, $mdNameA)
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrlA, $null, $null, $mdNameA)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $mdUrlB, $null, $null, $mdNameB)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdUrlA, $null, $null, $mdNameA)

$wsDe.Columns.Item(3).ColumnWidth = $wideTarget
$wsDe.Columns.Item(9).ColumnWidth = $fortyTarget
$wsDe.Columns.Item(10).ColumnWidth = $fortyTarget
